# Auto update stock data
# Updates the "as of" date (column A) and EBITDA (column B) for each
# company's latest-data row (2025/11/14 -> 2025/11/15), along with the
# corresponding EBITDA value. A leading apostrophe is used so Excel keeps
# these values as plain text (matching the source data, which stores
# dates/numbers as text) instead of auto-converting the date-like /
# numeric-like strings into a real date or number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 1).Value = "'2025/11/15"
}

$ebitda = @{
    2  = "4.58"
    8  = "7.60"
    14 = "2.84"
    20 = "12.00"
    26 = "9.68"
    32 = "24.62"
    44 = "10.49"
    50 = "11.22"
    56 = "32.48"
    68 = "12.53"
    74 = "15.04"
}

foreach ($r in $ebitda.Keys) {
    $ws.Cells.Item($r, 2).Value = "'" + $ebitda[$r]
}
